# Data Dictionary.xlsx - "Added all the relevant files" commit
#
# The variable that used to be documented as "avg_povrate" is renamed to
# "avg_povpoprate", and its description is updated to reflect that the
# poverty rate is now weighted by population size.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A48").Value = "avg_povpoprate"
$ws.Range("B48").Value = 'Average poverty rate weighted bypopulation size of the variable "zips" between 2012 and 2017'

# Move the selection/cursor down to where the sheet was left (row 51),
# matching the updated view state saved with the workbook.
$ws.Range("B51").Select()
